# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col D) and "Correspond Handback
# DateTime" (col G) for the 84b09259... row (row 2) on both the "zh-cn" and
# "de-de" language report sheets, reflecting the freshly generated handback
# report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-26 09:06:57"
$wsZhCn.Range("G2").Value = "2016-01-26 09:07:41"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-26 09:07:08"
$wsDeDe.Range("G2").Value = "2016-01-26 09:07:59"
